$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New attendance row for Afaq (2024-08-15 check-in) appended below the
# existing data (row 4 was the previous last row).
$ws.Range("A5").Value = "Afaq"

# "2024-08-15" looks like a date, and a plain .Value assignment would get
# silently parsed/stored as a date serial instead of literal text (the
# workbook stores it as shared text everywhere else). Enter it as a text
# formula first, then convert that formula to its static text result so
# the stored cell is plain shared-string text with no special formatting.
$ws.Range("B5").Formula = '="2024-08-15"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)

$ws.Range("C5").Value = "06:29:24 PM"
